$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")
$tbl = $ws.ListObjects.Item("Tabell3")

# New row 1 (sheet row 7)
$null = $tbl.ListRows.Add()
$ws.Range("A7").Value = 45310
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Spara aktivitet"

# New row 2 (sheet row 8)
$null = $tbl.ListRows.Add()
$ws.Range("A8").Value = 45310
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Hämta enskild aktivitet"

# New row 3 (sheet row 9)
$null = $tbl.ListRows.Add()
$ws.Range("A9").Value = 45310
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "Uppdatera aktivitet"

$ws.Activate()
$null = $ws.Range("B12").Select()
